function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. Update "总计" (summary) sheet: insert a new row for 2022-Q3
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).ClearFormats()
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 20
$ws1.Range("D2").Value = 26.98
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it
#    inherits the exact same column layout / header / styling), then
#    move it into place right after "总计" and rename it.
# -----------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$new = $wb.Worksheets.Item("2022-Q2 (2)")
$new.Name = "2022-Q3"

# The new sheet needs 20 data rows (rows 2-21); the template only had
# 15 (rows 2-16), so extend it with 5 more rows, matching the A-column
# style (bold/centered/bordered) used by the existing index cells.
$new.Range("A17").Value = 15
$new.Range("A18").Value = 16
$new.Range("A19").Value = 17
$new.Range("A20").Value = 18
$new.Range("A21").Value = 19
$new.Range("A16").Copy()
$new.Range("A17:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false


# Overwrite all data rows (2-21) with the 2022-Q3 figures.
$new.Range("A2").Value = 0
Set-TextValue $new.Range("B2") "004997"
Set-TextValue $new.Range("C2") "广发高端制造股票A"
Set-TextValue $new.Range("D2") "130.03"
Set-TextValue $new.Range("E2") "87.00"
Set-TextValue $new.Range("F2") "5.29"
Set-TextValue $new.Range("G2") "6.8786"
$new.Range("H2").Value = 8
$new.Range("A3").Value = 1
Set-TextValue $new.Range("B3") "002939"
Set-TextValue $new.Range("C3") "广发创新升级灵活配置混合"
Set-TextValue $new.Range("D3") "84.58"
Set-TextValue $new.Range("E3") "94.36"
Set-TextValue $new.Range("F3") "6.20"
Set-TextValue $new.Range("G3") "5.2440"
$new.Range("H3").Value = 9
$new.Range("A4").Value = 2
Set-TextValue $new.Range("B4") "162703"
Set-TextValue $new.Range("C4") "广发小盘成长混合（LOF）A"
Set-TextValue $new.Range("D4") "99.59"
Set-TextValue $new.Range("E4") "88.79"
Set-TextValue $new.Range("F4") "4.94"
Set-TextValue $new.Range("G4") "4.9197"
$new.Range("H4").Value = 8
$new.Range("A5").Value = 3
Set-TextValue $new.Range("B5") "002132"
Set-TextValue $new.Range("C5") "广发鑫享灵活配置混合A"
Set-TextValue $new.Range("D5") "65.56"
Set-TextValue $new.Range("E5") "88.72"
Set-TextValue $new.Range("F5") "3.72"
Set-TextValue $new.Range("G5") "2.4388"
$new.Range("H5").Value = 9
$new.Range("A6").Value = 4
Set-TextValue $new.Range("B6") "011130"
Set-TextValue $new.Range("C6") "广发兴诚混合C"
Set-TextValue $new.Range("D6") "25.57"
Set-TextValue $new.Range("E6") "88.20"
Set-TextValue $new.Range("F6") "7.22"
Set-TextValue $new.Range("G6") "1.8462"
$new.Range("H6").Value = 6
$new.Range("A7").Value = 5
Set-TextValue $new.Range("B7") "011479"
Set-TextValue $new.Range("C7") "广发诚享混合A"
Set-TextValue $new.Range("D7") "36.88"
Set-TextValue $new.Range("E7") "87.12"
Set-TextValue $new.Range("F7") "4.36"
Set-TextValue $new.Range("G7") "1.6080"
$new.Range("H7").Value = 10
$new.Range("A8").Value = 6
Set-TextValue $new.Range("B8") "011121"
Set-TextValue $new.Range("C8") "广发兴诚混合A"
Set-TextValue $new.Range("D8") "22.16"
Set-TextValue $new.Range("E8") "88.20"
Set-TextValue $new.Range("F8") "7.22"
Set-TextValue $new.Range("G8") "1.6000"
$new.Range("H8").Value = 6
$new.Range("A9").Value = 7
Set-TextValue $new.Range("B9") "013141"
Set-TextValue $new.Range("C9") "中信保诚弘远混合A"
Set-TextValue $new.Range("D9") "19.17"
Set-TextValue $new.Range("E9") "79.40"
Set-TextValue $new.Range("F9") "3.35"
Set-TextValue $new.Range("G9") "0.6422"
$new.Range("H9").Value = 6
$new.Range("A10").Value = 8
Set-TextValue $new.Range("B10") "165516"
Set-TextValue $new.Range("C10") "信诚周期轮动混合（LOF）A"
Set-TextValue $new.Range("D10") "17.70"
Set-TextValue $new.Range("E10") "82.67"
Set-TextValue $new.Range("F10") "3.28"
Set-TextValue $new.Range("G10") "0.5806"
$new.Range("H10").Value = 6
$new.Range("A11").Value = 9
Set-TextValue $new.Range("B11") "010160"
Set-TextValue $new.Range("C11") "广发高端制造股票C"
Set-TextValue $new.Range("D11") "7.66"
Set-TextValue $new.Range("E11") "87.00"
Set-TextValue $new.Range("F11") "5.29"
Set-TextValue $new.Range("G11") "0.4052"
$new.Range("H11").Value = 8
$new.Range("A12").Value = 10
Set-TextValue $new.Range("B12") "009132"
Set-TextValue $new.Range("C12") "广发小盘成长混合（LOF）C"
Set-TextValue $new.Range("D12") "5.90"
Set-TextValue $new.Range("E12") "88.79"
Set-TextValue $new.Range("F12") "4.94"
Set-TextValue $new.Range("G12") "0.2915"
$new.Range("H12").Value = 8
$new.Range("A13").Value = 11
Set-TextValue $new.Range("B13") "015322"
Set-TextValue $new.Range("C13") "广发鑫享灵活配置混合C"
Set-TextValue $new.Range("D13") "5.01"
Set-TextValue $new.Range("E13") "88.72"
Set-TextValue $new.Range("F13") "3.72"
Set-TextValue $new.Range("G13") "0.1864"
$new.Range("H13").Value = 9
$new.Range("A14").Value = 12
Set-TextValue $new.Range("B14") "011480"
Set-TextValue $new.Range("C14") "广发诚享混合C"
Set-TextValue $new.Range("D14") "4.20"
Set-TextValue $new.Range("E14") "87.12"
Set-TextValue $new.Range("F14") "4.36"
Set-TextValue $new.Range("G14") "0.1831"
$new.Range("H14").Value = 10
$new.Range("A15").Value = 13
Set-TextValue $new.Range("B15") "007074"
Set-TextValue $new.Range("C15") "国寿安保新蓝筹灵活配置混合"
Set-TextValue $new.Range("D15") "2.26"
Set-TextValue $new.Range("E15") "88.44"
Set-TextValue $new.Range("F15") "2.89"
Set-TextValue $new.Range("G15") "0.0653"
$new.Range("H15").Value = 9
$new.Range("A16").Value = 14
Set-TextValue $new.Range("B16") "001742"
Set-TextValue $new.Range("C16") "广发百发大数据策略精选灵活配置混合E"
Set-TextValue $new.Range("D16") "2.51"
Set-TextValue $new.Range("E16") "40.85"
Set-TextValue $new.Range("F16") "2.22"
Set-TextValue $new.Range("G16") "0.0557"
$new.Range("H16").Value = 9
$new.Range("A17").Value = 15
Set-TextValue $new.Range("B17") "010821"
Set-TextValue $new.Range("C17") "东方红多元策略混合B"
Set-TextValue $new.Range("D17") "0.59"
Set-TextValue $new.Range("E17") "92.90"
Set-TextValue $new.Range("F17") "3.09"
Set-TextValue $new.Range("G17") "0.0182"
$new.Range("H17").Value = 10
$new.Range("A18").Value = 16
Set-TextValue $new.Range("B18") "910017"
Set-TextValue $new.Range("C18") "东方红多元策略混合A"
Set-TextValue $new.Range("D18") "0.41"
Set-TextValue $new.Range("E18") "92.90"
Set-TextValue $new.Range("F18") "3.09"
Set-TextValue $new.Range("G18") "0.0127"
$new.Range("H18").Value = 10
$new.Range("A19").Value = 17
Set-TextValue $new.Range("B19") "001741"
Set-TextValue $new.Range("C19") "广发百发大数据策略精选灵活配置混合A"
Set-TextValue $new.Range("D19") "0.21"
Set-TextValue $new.Range("E19") "40.85"
Set-TextValue $new.Range("F19") "2.22"
Set-TextValue $new.Range("G19") "0.0047"
$new.Range("H19").Value = 9
$new.Range("A20").Value = 18
Set-TextValue $new.Range("B20") "014335"
Set-TextValue $new.Range("C20") "信诚周期轮动混合（LOF）C"
Set-TextValue $new.Range("D20") "0.03"
Set-TextValue $new.Range("E20") "82.67"
Set-TextValue $new.Range("F20") "3.28"
Set-TextValue $new.Range("G20") "0.0010"
$new.Range("H20").Value = 6
$new.Range("A21").Value = 19
Set-TextValue $new.Range("B21") "015936"
Set-TextValue $new.Range("C21") "中信保诚弘远混合C"
Set-TextValue $new.Range("D21") "0.02"
Set-TextValue $new.Range("E21") "79.40"
Set-TextValue $new.Range("F21") "3.35"
Set-TextValue $new.Range("G21") "0.0007"
$new.Range("H21").Value = 6
